# Add a new "sampleName" column (K) that concatenates Subject (A),
# sampleType (E) and sampleDate2 (J) for every data row, e.g.
# "197 ETA 20200330".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("K1").Value = "sampleName"

# Last row with data (header is row 1, data runs through row 50)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 50 }

for ($row = 2; $row -le $lastRow; $row++) {
    $subject = $ws.Cells.Item($row, 1).Value2
    $sampleType = $ws.Cells.Item($row, 5).Value2
    $sampleDate2 = $ws.Cells.Item($row, 10).Value2

    $ws.Range("K" + $row).Value = "$subject $sampleType $sampleDate2"
}

Write-Host "Added sampleName column (K) for rows 2..$lastRow"
